$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "58-25="
$t.Cell(1,2).Range.Text = "36-3="
$t.Cell(1,3).Range.Text = "10+25="
$t.Cell(1,4).Range.Text = "26-9="
$t.Cell(1,5).Range.Text = "55-7="
$t.Cell(2,1).Range.Text = "72-31="
$t.Cell(2,2).Range.Text = "9+64="
$t.Cell(2,3).Range.Text = "55+10="
$t.Cell(2,4).Range.Text = "63+0="
$t.Cell(2,5).Range.Text = "53+38="
$t.Cell(3,1).Range.Text = "46+16="
$t.Cell(3,2).Range.Text = "19+55="
$t.Cell(3,3).Range.Text = "62-4="
$t.Cell(3,4).Range.Text = "57+21="
$t.Cell(3,5).Range.Text = "56+11="
$t.Cell(4,1).Range.Text = "94-78="
$t.Cell(4,2).Range.Text = "80+17="
$t.Cell(4,3).Range.Text = "8+48="
$t.Cell(4,4).Range.Text = "44+44="
$t.Cell(4,5).Range.Text = "89-39="
$t.Cell(5,1).Range.Text = "63+16="
$t.Cell(5,2).Range.Text = "52+14="
$t.Cell(5,3).Range.Text = "28-11="
$t.Cell(5,4).Range.Text = "84-33="
$t.Cell(5,5).Range.Text = "27+32="
$t.Cell(6,1).Range.Text = "77+19="
$t.Cell(6,2).Range.Text = "12+87="
$t.Cell(6,3).Range.Text = "87-70="
$t.Cell(6,4).Range.Text = "43-6="
$t.Cell(6,5).Range.Text = "3+13="
$t.Cell(7,1).Range.Text = "0+0="
$t.Cell(7,2).Range.Text = "25+59="
$t.Cell(7,3).Range.Text = "73-60="
$t.Cell(7,4).Range.Text = "66-58="
$t.Cell(7,5).Range.Text = "37+58="
$t.Cell(8,1).Range.Text = "9+7="
$t.Cell(8,2).Range.Text = "4-2="
$t.Cell(8,3).Range.Text = "11+15="
$t.Cell(8,4).Range.Text = "55-22="
$t.Cell(8,5).Range.Text = "28+32="
$t.Cell(9,1).Range.Text = "4+49="
$t.Cell(9,2).Range.Text = "7+90="
$t.Cell(9,3).Range.Text = "93-27="
$t.Cell(9,4).Range.Text = "85-31="
$t.Cell(9,5).Range.Text = "27+48="
$t.Cell(10,1).Range.Text = "70-19="
$t.Cell(10,2).Range.Text = "27+15="
$t.Cell(10,3).Range.Text = "21+6="
$t.Cell(10,4).Range.Text = "3+40="
$t.Cell(10,5).Range.Text = "65-58="
$t.Cell(11,1).Range.Text = "77-74="
$t.Cell(11,2).Range.Text = "95-56="
$t.Cell(11,3).Range.Text = "44-28="
$t.Cell(11,4).Range.Text = "58+34="
$t.Cell(11,5).Range.Text = "97-85="
$t.Cell(12,1).Range.Text = "75-56="
$t.Cell(12,2).Range.Text = "86-81="
$t.Cell(12,3).Range.Text = "36-5="
$t.Cell(12,4).Range.Text = "73-11="
$t.Cell(12,5).Range.Text = "46+33="
$t.Cell(13,1).Range.Text = "20+24="
$t.Cell(13,2).Range.Text = "22-9="
$t.Cell(13,3).Range.Text = "8+75="
$t.Cell(13,4).Range.Text = "17+66="
$t.Cell(13,5).Range.Text = "55+35="
$t.Cell(14,1).Range.Text = "7+43="
$t.Cell(14,2).Range.Text = "0+19="
$t.Cell(14,3).Range.Text = "88-22="
$t.Cell(14,4).Range.Text = "0+34="
$t.Cell(14,5).Range.Text = "91+6="
$t.Cell(15,1).Range.Text = "31-7="
$t.Cell(15,2).Range.Text = "19+46="
$t.Cell(15,3).Range.Text = "1+62="
$t.Cell(15,4).Range.Text = "38-5="
$t.Cell(15,5).Range.Text = "82-35="
$t.Cell(16,1).Range.Text = "10+47="
$t.Cell(16,2).Range.Text = "47-31="
$t.Cell(16,3).Range.Text = "25+0="
$t.Cell(16,4).Range.Text = "1+18="
$t.Cell(16,5).Range.Text = "25+69="
$t.Cell(17,1).Range.Text = "62-4="
$t.Cell(17,2).Range.Text = "16+5="
$t.Cell(17,3).Range.Text = "54-16="
$t.Cell(17,4).Range.Text = "92-34="
$t.Cell(17,5).Range.Text = "94-19="
$t.Cell(18,1).Range.Text = "2+54="
$t.Cell(18,2).Range.Text = "65+34="
$t.Cell(18,3).Range.Text = "99-91="
$t.Cell(18,4).Range.Text = "99-44="
$t.Cell(18,5).Range.Text = "12+42="
$t.Cell(19,1).Range.Text = "29+23="
$t.Cell(19,2).Range.Text = "83-18="
$t.Cell(19,3).Range.Text = "77-42="
$t.Cell(19,4).Range.Text = "45-23="
$t.Cell(19,5).Range.Text = "35-17="
$t.Cell(20,1).Range.Text = "66-60="
$t.Cell(20,2).Range.Text = "11+9="
$t.Cell(20,3).Range.Text = "23-16="
$t.Cell(20,4).Range.Text = "66+19="
$t.Cell(20,5).Range.Text = "69-64="
